$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# Row 2 (Cyprus, company id "3")
# ----------------------------------------------------------------------------
# B2 holds a numeric-looking identifier that must stay text (matches the
# original inline string "2" -> "3"). Force text formatting so Excel does
# not silently coerce it to a number, then drop the per-cell format
# override again so the cell keeps the default (unstyled) appearance.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("B2").Style = "Normal"

$ws.Range("D2").Value = 0.218
$ws.Range("E2").Value = 0.6595
$ws.Range("F2").Value = 0.191
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 260.8000000000001
$ws.Range("L2").Value = 0.1348988775668546
$ws.Range("M2").Value = 112.31
$ws.Range("N2").Value = 0.01550643397581046
$ws.Range("O2").Value = 0.4306365030674846
$ws.Range("P2").Value = 103.8
$ws.Range("Q2").Value = 0.01433147401557409
$ws.Range("R2").Value = 0.398006134969325
$ws.Range("S2").Value = 8.510000000000005
$ws.Range("T2").Value = 0.07577241563529521
$ws.Range("U2").Value = 10432.5
$ws.Range("V2").Value = 1.440395979455459
$ws.Range("W2").Value = 0.06024409244352116
$ws.Range("X2").Value = 0.05856591019418327
$ws.Range("Y2").Value = 0.001678182249337898
$ws.Range("Z2").Value = -0.5396812103955561
$ws.Range("AA2").Value = -0
$ws.Range("AB2").Value = 0.05329997332913817
$ws.Range("AC2").Value = -0.05329997332913817
$ws.Range("AD2").Value = 2550.7
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2550.7
$ws.Range("AG2").Value = -7881.8
$ws.Range("AH2").Value = 0.2604482564966559
$ws.Range("AI2").Value = 0.3244876410497793
$ws.Range("AJ2").Value = 12.33458528951487
$ws.Range("AK2").Value = 3.064701765300568
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# ----------------------------------------------------------------------------
# Row 3 (Cyprus, TCS Group Holding PLC (LSE:TCS))
# ----------------------------------------------------------------------------
$ws.Range("B3").Value = "TCS Group Holding PLC (LSE:TCS)"
$ws.Range("D3").Value = 0.535
$ws.Range("E3").Value = 0.95
$ws.Range("F3").Value = 0.191
$ws.Range("K3").Value = 552.2
$ws.Range("L3").Value = 0.472693032015066
$ws.Range("M3").Value = 112.31
$ws.Range("N3").Value = 0.0173908330752555
$ws.Range("O3").Value = 0.2033864541832669
$ws.Range("P3").Value = 103.8
$ws.Range("Q3").Value = 0.0160730876432332
$ws.Range("R3").Value = 0.1879753712423035
$ws.Range("S3").Value = 8.510000000000005
$ws.Range("T3").Value = 0.07577241563529521
$ws.Range("U3").Value = 636.4
$ws.Range("V3").Value = 0.09854444100340662
$ws.Range("W3").Value = 0.4306324573032832
$ws.Range("X3").Value = 0.04938574416024047
$ws.Range("Y3").Value = 0.3812467131430427
$ws.Range("Z3").Value = 0.7154142935880948
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04845915785033052
$ws.Range("AC3").Value = -0.04845915785033052
$ws.Range("AD3").Value = 709.9
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 709.9
$ws.Range("AG3").Value = 73.5
$ws.Range("AH3").Value = 0.09903877007212712
$ws.Range("AI3").Value = 0.3212071851952401
$ws.Range("AJ3").Value = 0.01125315777386512
$ws.Range("AK3").Value = 0.04670521700451166
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0

# ----------------------------------------------------------------------------
# Row 4 (new row: Cyprus, Hellenic Bank Public Company Limited (CSE:HB))
# ----------------------------------------------------------------------------
$ws.Range("A4").Value = "Cyprus"
$ws.Range("B4").Value = "Hellenic Bank Public Company Limited (CSE:HB)"
$ws.Range("C4").Value = "Bank (Money Center)"
$ws.Range("D4").Value = 0.207
$ws.Range("E4").Value = 0.369
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 69.59999999999999
$ws.Range("L4").Value = 0.1592677345537757
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("S4").Value = 0
$ws.Range("T4").ClearContents()
$ws.Range("U4").Value = 3330.2
$ws.Range("V4").Value = 8.654365904365903
$ws.Range("W4").Value = 0.06024409244352116
$ws.Range("X4").Value = 0.05856591019418327
$ws.Range("Y4").Value = 0.001678182249337898
$ws.Range("Z4").Value = -0.1176850779629979
$ws.Range("AA4").Value = -0
$ws.Range("AB4").Value = 0.05329997332913817
$ws.Range("AC4").Value = -0.05329997332913817
$ws.Range("AD4").Value = 152.7
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 152.7
$ws.Range("AG4").Value = -3177.5
$ws.Range("AH4").Value = 0.2840930232558139
$ws.Range("AI4").Value = 0.1045317634173056
$ws.Range("AJ4").Value = 1.137787803917356
$ws.Range("AK4").Value = 1.699743233122927
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# ----------------------------------------------------------------------------
# Row 5 (new row: Cyprus, Bank of Cyprus Holdings Public Limited Company (CSE:BOCH))
# ----------------------------------------------------------------------------
$ws.Range("A5").Value = "Cyprus"
$ws.Range("B5").Value = "Bank of Cyprus Holdings Public Limited Company (CSE:BOCH)"
$ws.Range("C5").Value = "Bank (Money Center)"
$ws.Range("D5").Value = 0.218
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -361
$ws.Range("L5").Value = -1.100274306613837
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 6465.9
$ws.Range("V5").Value = 16.16475
$ws.Range("W5").Value = -0.1348977990359105
$ws.Range("X5").Value = 0.1809052758796562
$ws.Range("Y5").Value = -0.3158030749155667
$ws.Range("Z5").Value = -0.2184566216126241
$ws.Range("AA5").Value = -0
$ws.Range("AB5").Value = 0.0716849993168773
$ws.Range("AC5").Value = -0.0716849993168773
$ws.Range("AD5").Value = 1688.1
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 1688.1
$ws.Range("AG5").Value = -4777.799999999999
$ws.Range("AH5").Value = 0.8084382931851922
$ws.Range("AI5").Value = 0.4029070600028641
$ws.Range("AJ5").Value = 1.091370094568048
$ws.Range("AK5").Value = 2.099116910504811
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0

Write-Output "done"
